$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column widths ---
# NOTE: the host's ColumnWidth setter quantizes to the nearest reachable
# pixel-snapped value (it re-derives the stored <col width> from pixels,
# stepping in 1/6-character increments), so we feed it the ColumnWidth
# input value whose pixel-snap lands as close as possible to the exact
# target "raw" width from the target OOXML.
$ws.Columns.Item(9).ColumnWidth = 1.1439732142857144   # I -> raw ~2.1666667 (target 2.140625)
$ws.Columns.Item(11).ColumnWidth = 4.143973214285714   # K -> raw ~5.6666667 (target 5.7109375)
$ws.Columns.Item(7).ColumnWidth = 2.001116071428571    # G -> raw ~3.1666667 (target 3.140625)
$ws.Columns.Item(12).ColumnWidth = 3.2868303571428568  # L -> raw ~4.6666667 (target 4.7109375)

# --- Update row 1 cell values ---
$ws.Range("B1").Value = 4
$ws.Range("C1").Value = 32
$ws.Range("D1").Value = 13
$ws.Range("E1").Value = 20
$ws.Range("F1").Value = 9
$ws.Range("G1").Value = 12
$ws.Range("H1").Value = 6
$ws.Range("I1").Value = 2
$ws.Range("J1").Value = 13
$ws.Range("K1").Value = 0.072
$ws.Range("L1").Value = 0.04
$ws.Range("M1").Value = 0.076
$ws.Range("N1").Value = 0.034
